# Updated Hybrid framework with all runmode set to N
$wb = $excel.ActiveWorkbook

# TestCases sheet: flip Runmode column (B2:B4) from "Y" to "N"
$ws1 = $wb.Worksheets.Item("TestCases")
$ws1.Range("B2:B4").Value = "N"

# Make TestCases the active sheet/tab, with B2:B4 (the edited range) selected
[void]$ws1.Activate()
[void]$ws1.Range("B2:B4").Select()
